# Update the "as_of_utc" timestamp column (AA) on the "Главные" and
# "Линейные" worksheets from "2025-11-23 03:05:04" to "2025-11-23 07:05:06"
# for every data row (rows 2-26).

$wb = $excel.ActiveWorkbook

$oldValue = "2025-11-23 03:05:04"
$newValue = "2025-11-23 07:05:06"

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $cell = $ws.Range("AA$row")
        if ($cell.Value() -eq $oldValue) {
            $cell.Value = $newValue
        }
    }
}
